# Fix typo in cell B10: "wilaNah" -> "wilayah"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "3.2 Peserta di luar`nwilayah kerja"
